# Apply the 27-May-2024 GitHub Actions crypto-price refresh to sheet1.
# Source diff only ever touches columns B (coin name), C (link), D (price) and
# E (1h volume change) on existing rows; row order / A-column rank is untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.430.82'
$ws.Range("E2").Value = '  +0.79%  '

# Row 3
$ws.Range("D3").Value = '3.883.75'
$ws.Range("E3").Value = '  +0.36%  '

# Row 4
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.19%  '

# Row 5
$ws.Range("D5").Value = '''603.43'
$ws.Range("E5").Value = '  +0.42%  '

# Row 6
$ws.Range("D6").Value = '''169.32'
$ws.Range("E6").Value = '  +4.11%  '

# Row 7
$ws.Range("D7").Value = '3.884.33'
$ws.Range("E7").Value = '  +0.41%  '

# Row 8
$ws.Range("E8").Value = '  +0.07%  '

# Row 9
$ws.Range("D9").Value = '''0.534'
$ws.Range("E9").Value = '  +0.64%  '

# Row 10
$ws.Range("D10").Value = '''0.167'
$ws.Range("E10").Value = '  -0.50%  '

# Row 11
$ws.Range("E11").Value = '  +1.13%  '

# Row 12
$ws.Range("D12").Value = '''0.467'
$ws.Range("E12").Value = '  +1.84%  '

# Row 13
$ws.Range("D13").Value = '''0.0000255'
$ws.Range("E13").Value = '  +4.70%  '

# Row 14
$ws.Range("D14").Value = '''38.17'
$ws.Range("E14").Value = '  +3.33%  '

# Row 15
$ws.Range("D15").Value = '4.538.76'
$ws.Range("E15").Value = '  +0.52%  '

# Row 16
$ws.Range("D16").Value = '3.871.04'
$ws.Range("E16").Value = '  -0.62%  '

# Row 17
$ws.Range("D17").Value = '69.424.40'
$ws.Range("E17").Value = '  +0.54%  '

# Row 18
$ws.Range("D18").Value = '''18.74'
$ws.Range("E18").Value = '  +9.63%  '

# Row 19
$ws.Range("D19").Value = '''7.62'
$ws.Range("E19").Value = '  +1.08%  '

# Row 20
$ws.Range("E20").Value = '  -0.99%  '

# Row 21
$ws.Range("D21").Value = '''11.15'
$ws.Range("E21").Value = '  -1.46%  '

# Row 22
$ws.Range("D22").Value = '''489.95'
$ws.Range("E22").Value = '  +1.07%  '

# Row 23
$ws.Range("D23").Value = '''0.746'
$ws.Range("E23").Value = '  +3.71%  '

# Row 24
$ws.Range("D24").Value = '''0.0000166'
$ws.Range("E24").Value = '  +2.49%  '

# Row 25
$ws.Range("D25").Value = '''85.07'
$ws.Range("E25").Value = '  +1.25%  '

# Row 26
$ws.Range("E26").Value = '  +2.04%  '

# Row 27
$ws.Range("D27").Value = '''12.30'
$ws.Range("E27").Value = '  +1.77%  '

# Row 28
$ws.Range("D28").Value = '''10.14'
$ws.Range("E28").Value = '  +2.05%  '

# Row 29
$ws.Range("E29").Value = '  +0.10%  '

# Row 30
$ws.Range("E30").Value = '  +0.69%  '

# Row 31
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = '''2.42'
$ws.Range("E31").Value = '  +2.06%  '

# Row 32
$ws.Range("B32").Value = 'WrappedeETH'
$ws.Range("C32").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D32").Value = '4.034.23'
$ws.Range("E32").Value = '  +0.27%  '

# Row 33
$ws.Range("D33").Value = '''7.79'
$ws.Range("E33").Value = '  -1.44%  '

# Row 34
$ws.Range("D34").Value = '''31.87'
$ws.Range("E34").Value = '  -1.39%  '

# Row 35
$ws.Range("D35").Value = '3.847.62'
$ws.Range("E35").Value = '  +0.84%  '

# Row 36
$ws.Range("E36").Value = '  +0.18%  '

# Row 37
$ws.Range("D37").Value = '''6.10'
$ws.Range("E37").Value = '  +3.69%  '

# Row 39
$ws.Range("E39").Value = '  +0.80%  '

# Row 40
$ws.Range("D40").Value = '''3.28'
$ws.Range("E40").Value = '  +10.49%  '

# Row 41
$ws.Range("D41").Value = '''0.998'
$ws.Range("E41").Value = '  -0.18%  '

# Row 42
$ws.Range("E42").Value = '  +2.49%  '

# Row 43
$ws.Range("E43").Value = '  +6.35%  '

# Row 44
$ws.Range("D44").Value = '''436.17'
$ws.Range("E44").Value = '  -0.10%  '

# Row 45
$ws.Range("D45").Value = '''48.11'
$ws.Range("E45").Value = '  -0.79%  '

# Row 46
$ws.Range("D46").Value = '''8.69'
$ws.Range("E46").Value = '  +3.40%  '

# Row 47
$ws.Range("E47").Value = '  +0.02%  '

# Row 48
$ws.Range("D48").Value = '''0.0368'
$ws.Range("E48").Value = '  +3.27%  '

# Row 49
$ws.Range("D49").Value = '''143.55'
$ws.Range("E49").Value = '  +0.11%  '

# Row 50
$ws.Range("D50").Value = '''40.07'
$ws.Range("E50").Value = '  +3.34%  '

# Row 51
$ws.Range("D51").Value = '''0.000270'
$ws.Range("E51").Value = '  +18.80%  '
